# Re-order / relabel the header row (row 1) on every sheet of the
# 360-giving-schema-titles workbook, and add the two new required
# columns ("Activity/id:fundingOrganization" -> "Department" on
# Organization, "Funding Org:Department" on Activity). The Event sheet
# loses its "Activity/id:awardDate" column (merged into "Award Date" on
# Activity), so its header row shrinks by one column.

$wb = $excel.ActiveWorkbook

# --- Activity ---
$ws = $wb.Worksheets.Item("Activity")
$ws.Cells.Item(1, 1).Value = "Identifier"
$ws.Cells.Item(1, 2).Value = "Title"
$ws.Cells.Item(1, 3).Value = "Description"
$ws.Cells.Item(1, 4).Value = "Currency"
$ws.Cells.Item(1, 5).Value = "Amount Applied For"
$ws.Cells.Item(1, 6).Value = "Amount Awarded"
$ws.Cells.Item(1, 7).Value = "Amount Disbursed"
$ws.Cells.Item(1, 8).Value = "Award Date"
$ws.Cells.Item(1, 9).Value = "URL"
$ws.Cells.Item(1, 10).Value = "Planned Dates:Start Date"
$ws.Cells.Item(1, 11).Value = "Planned Dates:End Date"
$ws.Cells.Item(1, 12).Value = "Planned Dates:Duration (months)"
$ws.Cells.Item(1, 13).Value = "Recipient Org:Identifier"
$ws.Cells.Item(1, 14).Value = "Recipient Org:Name"
$ws.Cells.Item(1, 15).Value = "Recipient Org:Charity Number"
$ws.Cells.Item(1, 16).Value = "Recipient Org:Company Number"
$ws.Cells.Item(1, 17).Value = "Recipient Org:Street Address"
$ws.Cells.Item(1, 18).Value = "Recipient Org:City"
$ws.Cells.Item(1, 19).Value = "Recipient Org:Postal Code"
$ws.Cells.Item(1, 20).Value = "Recipient Org:Description"
$ws.Cells.Item(1, 21).Value = "Recipient Org:Web Address"
$ws.Cells.Item(1, 22).Value = "Beneficiary Location:Name"
$ws.Cells.Item(1, 23).Value = "Beneficiary Location:Country Code"
$ws.Cells.Item(1, 24).Value = "Beneficiary Location:Latitude"
$ws.Cells.Item(1, 25).Value = "Beneficiary Location:Longitude"
$ws.Cells.Item(1, 26).Value = "Beneficiary Location:Geographic Code"
$ws.Cells.Item(1, 27).Value = "Beneficiary Location:Geographic Code Type"
$ws.Cells.Item(1, 28).Value = "Funding Org:Identifier"
$ws.Cells.Item(1, 29).Value = "Funding Org:Name"
$ws.Cells.Item(1, 30).Value = "Funding Org:Department"
$ws.Cells.Item(1, 31).Value = "Grant Programme:Code"
$ws.Cells.Item(1, 32).Value = "Grant Programme:Title"
$ws.Cells.Item(1, 33).Value = "Grant Programme:URL"
$ws.Cells.Item(1, 34).Value = "From an open call?"
$ws.Cells.Item(1, 35).Value = "Related Activity"
$ws.Cells.Item(1, 36).Value = "Last modified"
$ws.Cells.Item(1, 37).Value = "Data Source"

# --- Classification ---
$ws = $wb.Worksheets.Item("Classification")
$ws.Cells.Item(1, 1).Value = "ocid"
$ws.Cells.Item(1, 2).Value = "Activity/id:fundingType"
$ws.Cells.Item(1, 3).Value = "Activity/id:classifications"
$ws.Cells.Item(1, 4).Value = "Vocabulary"
$ws.Cells.Item(1, 5).Value = "Code"
$ws.Cells.Item(1, 6).Value = "Title"
$ws.Cells.Item(1, 7).Value = "Description"
$ws.Cells.Item(1, 8).Value = "URL"
$ws.Cells.Item(1, 9).Value = "Last modified"

# --- Documents ---
$ws = $wb.Worksheets.Item("Documents")
$ws.Cells.Item(1, 1).Value = "ocid"
$ws.Cells.Item(1, 2).Value = "Activity/id:relatedDocument"
$ws.Cells.Item(1, 3).Value = "Identifier"
$ws.Cells.Item(1, 4).Value = "Title"
$ws.Cells.Item(1, 5).Value = "Web Address"
$ws.Cells.Item(1, 6).Value = "Description"
$ws.Cells.Item(1, 7).Value = "Document Type"
$ws.Cells.Item(1, 8).Value = "Last modified"

# --- Event ---
$ws = $wb.Worksheets.Item("Event")
$ws.Cells.Item(1, 1).Value = "ocid"
$ws.Cells.Item(1, 2).Value = "Activity/id:plannedDates"
$ws.Cells.Item(1, 3).Value = "Activity/id:actualDates"
$ws.Cells.Item(1, 4).Value = "Title"
$ws.Cells.Item(1, 5).Value = "Start Date"
$ws.Cells.Item(1, 6).Value = "End Date"
$ws.Cells.Item(1, 7).Value = "Duration (months)"
$ws.Cells.Item(1, 8).Value = "Description"
$ws.Cells.Item(1, 9).Value = "Last modified"
$ws.Cells.Item(1, 10).Clear()

# --- GrantProgramme ---
$ws = $wb.Worksheets.Item("GrantProgramme")
$ws.Cells.Item(1, 1).Value = "ocid"
$ws.Cells.Item(1, 2).Value = "Activity/id:grantProgramme"
$ws.Cells.Item(1, 3).Value = "Code"
$ws.Cells.Item(1, 4).Value = "Title"
$ws.Cells.Item(1, 5).Value = "Description"
$ws.Cells.Item(1, 6).Value = "URL"
$ws.Cells.Item(1, 7).Value = "Last modified"

# --- Location ---
$ws = $wb.Worksheets.Item("Location")
$ws.Cells.Item(1, 1).Value = "ocid"
$ws.Cells.Item(1, 2).Value = "Activity/id:location"
$ws.Cells.Item(1, 3).Value = "Activity/recipientOrganization[]/id:location"
$ws.Cells.Item(1, 4).Value = "Activity/id:beneficiaryLocation"
$ws.Cells.Item(1, 5).Value = "Activity/fundingOrganization[]/id:location"
$ws.Cells.Item(1, 6).Value = "Identifier"
$ws.Cells.Item(1, 7).Value = "Name"
$ws.Cells.Item(1, 8).Value = "Country Code"
$ws.Cells.Item(1, 9).Value = "Latitude"
$ws.Cells.Item(1, 10).Value = "Longitude"
$ws.Cells.Item(1, 11).Value = "Description"
$ws.Cells.Item(1, 12).Value = "Geographic Code"
$ws.Cells.Item(1, 13).Value = "Geographic Code Type"
$ws.Cells.Item(1, 14).Value = "Last modified"

# --- Organization ---
$ws = $wb.Worksheets.Item("Organization")
$ws.Cells.Item(1, 1).Value = "ocid"
$ws.Cells.Item(1, 2).Value = "Activity/id:recipientOrganization"
$ws.Cells.Item(1, 3).Value = "Activity/id:fundingOrganization"
$ws.Cells.Item(1, 4).Value = "Identifier"
$ws.Cells.Item(1, 5).Value = "Name"
$ws.Cells.Item(1, 6).Value = "Department"
$ws.Cells.Item(1, 7).Value = "Contact Name"
$ws.Cells.Item(1, 8).Value = "Charity Number"
$ws.Cells.Item(1, 9).Value = "Company Number"
$ws.Cells.Item(1, 10).Value = "Street Address"
$ws.Cells.Item(1, 11).Value = "City"
$ws.Cells.Item(1, 12).Value = "County"
$ws.Cells.Item(1, 13).Value = "Country"
$ws.Cells.Item(1, 14).Value = "Postal Code"
$ws.Cells.Item(1, 15).Value = "Phone Number"
$ws.Cells.Item(1, 16).Value = "Alternate Name"
$ws.Cells.Item(1, 17).Value = "Email"
$ws.Cells.Item(1, 18).Value = "Description"
$ws.Cells.Item(1, 19).Value = "Organisation Type"
$ws.Cells.Item(1, 20).Value = "Web Address"
$ws.Cells.Item(1, 21).Value = "Last modified"

# --- Transaction ---
$ws = $wb.Worksheets.Item("Transaction")
$ws.Cells.Item(1, 1).Value = "ocid"
$ws.Cells.Item(1, 2).Value = "Activity/id:applicationTransaction"
$ws.Cells.Item(1, 3).Value = "Activity/id:commitmentTransaction"
$ws.Cells.Item(1, 4).Value = "Activity/id:disbursementTransaction"
$ws.Cells.Item(1, 5).Value = "Identifier"
$ws.Cells.Item(1, 6).Value = "Transaction date"
$ws.Cells.Item(1, 7).Value = "Currency"
$ws.Cells.Item(1, 8).Value = "Value"
$ws.Cells.Item(1, 9).Value = "Value date"
$ws.Cells.Item(1, 10).Value = "Description"
$ws.Cells.Item(1, 11).Value = "Provider"
$ws.Cells.Item(1, 12).Value = "Recipient"
$ws.Cells.Item(1, 13).Value = "Last modified"
